$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") bumps from 46072 to 46073 for every data row (2-16)
foreach ($r in 2..16) {
    $ws.Cells.Item($r, 3).Value = 46073
}

# Rows 8, 9, 11, 12, 14, 15 are re-synced with new Beteckning/Datum/Area values
# (rows 10, 13, 16 keep their existing A/B/G values, only C changes above).

# Row 8: was "A 58926-2025", now "A 50530-2024"
$ws.Cells.Item(8, 1).Value = "A 50530-2024"
$ws.Cells.Item(8, 2).Value = 45601.56424768519
$ws.Cells.Item(8, 7).Value = 0.7

# Row 9: was "A 45370-2022", now "A 58926-2025"
$ws.Cells.Item(9, 1).Value = "A 58926-2025"
$ws.Cells.Item(9, 2).Value = 45986
$ws.Cells.Item(9, 7).Value = 3.1

# Row 11: was "A 2253-2022", now "A 23678-2023"
$ws.Cells.Item(11, 1).Value = "A 23678-2023"
$ws.Cells.Item(11, 2).Value = 45077
$ws.Cells.Item(11, 7).Value = 1.4

# Row 12: was "A 23678-2023", now "A 50538-2024"
$ws.Cells.Item(12, 1).Value = "A 50538-2024"
$ws.Cells.Item(12, 2).Value = 45601.57153935185
$ws.Cells.Item(12, 7).Value = 0.8

# Row 14: was "A 50530-2024", now "A 45370-2022"
$ws.Cells.Item(14, 1).Value = "A 45370-2022"
$ws.Cells.Item(14, 2).Value = 44844.6397337963
$ws.Cells.Item(14, 7).Value = 2.7

# Row 15: was "A 50538-2024", now "A 2253-2022"
$ws.Cells.Item(15, 1).Value = "A 2253-2022"
$ws.Cells.Item(15, 2).Value = 44578
$ws.Cells.Item(15, 7).Value = 0.3
